$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3473287.8
$ws.Range("I132").Value = 1099.871
$ws.Range("K132").Value = 3299.613
$ws.Range("M132").Value = -769.6130000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6281783.5
$ws.Range("I61").Value = 3206231.2
$ws.Range("J61").Value = 19609176
$ws.Range("K61").Value = 3206231.2
$ws.Range("L61").Value = 19609176
$ws.Range("M61").Value = -3206019.2
$ws.Range("N61").Value = -19609600
$ws.Range("H122").Value = 1742.0834
$ws.Range("I122").Value = 1157.9286
$ws.Range("J122").Value = 2559.9
$ws.Range("K122").Value = 3473.7858
$ws.Range("L122").Value = 7679.700000000001
$ws.Range("M122").Value = -1023.7858
$ws.Range("N122").Value = -12579.7
$ws.Range("H136").Value = 6281783.5
$ws.Range("I136").Value = 3206231.2
$ws.Range("J136").Value = 19609176
$ws.Range("K136").Value = 9618693.600000001
$ws.Range("L136").Value = 58827528
$ws.Range("M136").Value = -9616143.600000001
$ws.Range("N136").Value = -58832628
$ws.Range("H139").Value = 32771.715
$ws.Range("J139").Value = 32771.715
$ws.Range("L139").Value = 32771.715
$ws.Range("N139").Value = -43051.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17858460
$ws.Range("I134").Value = 23810896
$ws.Range("J134").Value = 3969438.8
$ws.Range("K134").Value = 71432688
$ws.Range("L134").Value = 11908316.4
$ws.Range("M134").Value = -71430153
$ws.Range("N134").Value = -11913386.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4533146
$ws.Range("I31").Value = 2977174.5
$ws.Range("J31").Value = 6953546
$ws.Range("K31").Value = 2977174.5
$ws.Range("L31").Value = 6953546
$ws.Range("M31").Value = -2976879.5
$ws.Range("N31").Value = -6954136
$ws.Range("H34").Value = 4533146
$ws.Range("I34").Value = 2977174.5
$ws.Range("J34").Value = 6953546
$ws.Range("K34").Value = 2977174.5
$ws.Range("L34").Value = 6953546
$ws.Range("M34").Value = -2976972.5
$ws.Range("N34").Value = -6953950
$ws.Range("H58").Value = 4775850.5
$ws.Range("I58").Value = 3247810.2
$ws.Range("K58").Value = 3247810.2
$ws.Range("M58").Value = -3247607.2
$ws.Range("H80").Value = 32500
$ws.Range("I80").Value = 30000
$ws.Range("J80").Value = 35000
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 35000
$ws.Range("M80").Value = -28877
$ws.Range("N80").Value = -37246
$ws.Range("H83").Value = 32500
$ws.Range("I83").Value = 30000
$ws.Range("J83").Value = 35000
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 105000
$ws.Range("M83").Value = -84384
$ws.Range("N83").Value = -116232
$ws.Range("H132").Value = 3847927.5
$ws.Range("I132").Value = 6250881
$ws.Range("J132").Value = 3202.4
$ws.Range("K132").Value = 18752643
$ws.Range("L132").Value = 9607.200000000001
$ws.Range("M132").Value = -18750113
$ws.Range("N132").Value = -14667.2
$ws.Range("H134").Value = 1291735.9
$ws.Range("I134").Value = 1327.5555
$ws.Range("J134").Value = 3078455.2
$ws.Range("K134").Value = 3982.6665
$ws.Range("L134").Value = 9235365.600000001
$ws.Range("M134").Value = -1447.6665
$ws.Range("N134").Value = -9240435.600000001
$ws.Range("H136").Value = 4775850.5
$ws.Range("I136").Value = 3247810.2
$ws.Range("K136").Value = 9743430.600000001
$ws.Range("M136").Value = -9740880.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2992336.2
$ws.Range("I5").Value = 2404521
$ws.Range("J5").Value = 4167967
$ws.Range("K5").Value = 7213563
$ws.Range("L5").Value = 12503901
$ws.Range("M5").Value = -7213451
$ws.Range("N5").Value = -12504125
$ws.Range("H122").Value = 518.34784
$ws.Range("I122").Value = 388.9375
$ws.Range("J122").Value = 814.1429000000001
$ws.Range("K122").Value = 3500.4375
$ws.Range("L122").Value = 7327.2861
$ws.Range("M122").Value = -1050.4375
$ws.Range("N122").Value = -12227.2861
$ws.Range("H135").Value = 2992336.2
$ws.Range("I135").Value = 2404521
$ws.Range("J135").Value = 4167967
$ws.Range("K135").Value = 21640689
$ws.Range("L135").Value = 37511703
$ws.Range("M135").Value = -21638154
$ws.Range("N135").Value = -37516773
$ws.Range("H139").Value = 72214.42999999999
$ws.Range("I139").Value = 83916.836
$ws.Range("K139").Value = 251750.508
$ws.Range("M139").Value = -246610.508

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19500
$ws.Range("J15").Value = 19500
$ws.Range("L15").Value = 19500
$ws.Range("N15").Value = -20076
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H81").Value = 19500
$ws.Range("J81").Value = 19500
$ws.Range("L81").Value = 19500
$ws.Range("N81").Value = -21496
$ws.Range("H84").Value = 19500
$ws.Range("J84").Value = 19500
$ws.Range("L84").Value = 58500
$ws.Range("N84").Value = -68484
$ws.Range("H126").Value = 14233
$ws.Range("I126").Value = 15923.429
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 47770.287
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -45300.287
$ws.Range("N126").Value = -12140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1772.2858
$ws.Range("I7").Value = 1681.6
$ws.Range("K7").Value = 1681.6
$ws.Range("M7").Value = -1569.6
$ws.Range("H61").Value = 2037.4
$ws.Range("I61").Value = 1400.5
$ws.Range("J61").Value = 2462
$ws.Range("K61").Value = 1400.5
$ws.Range("L61").Value = 2462
$ws.Range("M61").Value = -1198.5
$ws.Range("N61").Value = -2866
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 2037.4
$ws.Range("I113").Value = 1400.5
$ws.Range("J113").Value = 2462
$ws.Range("K113").Value = 1400.5
$ws.Range("L113").Value = 2462
$ws.Range("M113").Value = 769.5
$ws.Range("N113").Value = -6802
$ws.Range("H126").Value = 1772.2858
$ws.Range("I126").Value = 1681.6
$ws.Range("K126").Value = 5044.799999999999
$ws.Range("M126").Value = -2574.799999999999
$ws.Range("H132").Value = 1713781.8
$ws.Range("I132").Value = 2386075.2
$ws.Range("J132").Value = 2489.3635
$ws.Range("K132").Value = 7158225.600000001
$ws.Range("L132").Value = 7468.0905
$ws.Range("M132").Value = -7155695.600000001
$ws.Range("N132").Value = -12528.0905
$ws.Range("H136").Value = 2802662.5
$ws.Range("I136").Value = 3461702
$ws.Range("J136").Value = 1745
$ws.Range("K136").Value = 10385106
$ws.Range("L136").Value = 5235
$ws.Range("M136").Value = -10382556
$ws.Range("N136").Value = -10335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = 0
$ws.Range("H126").Value = 1776.3572
$ws.Range("I126").Value = 533.75
$ws.Range("J126").Value = 3433.1667
$ws.Range("K126").Value = 1601.25
$ws.Range("L126").Value = 10299.5001
$ws.Range("M126").Value = 868.75
$ws.Range("N126").Value = -15239.5001
$ws.Range("H132").Value = 2027408.6
$ws.Range("I132").Value = 1818570.9
$ws.Range("J132").Value = 2332633
$ws.Range("K132").Value = 5455712.699999999
$ws.Range("L132").Value = 6997899
$ws.Range("M132").Value = -5453182.699999999
$ws.Range("N132").Value = -7002959
$ws.Range("H133").Value = 52134.5
$ws.Range("J133").Value = 52134.5
$ws.Range("L133").Value = 52134.5
$ws.Range("N133").Value = -62254.5
$ws.Range("H136").Value = 1793.5135
$ws.Range("I136").Value = 677.3
$ws.Range("J136").Value = 3106.7058
$ws.Range("K136").Value = 2031.9
$ws.Range("L136").Value = 9320.117400000001
$ws.Range("M136").Value = 518.1000000000001
$ws.Range("N136").Value = -14420.1174
